$p = $ppt.ActivePresentation
Write-Host "Designs count before:" $p.Designs.Count
try {
  $d = $p.Designs.Add("MyDesign")
  Write-Host "Added design, count now:" $p.Designs.Count
} catch {
  Write-Host "ERROR ADD: $_"
}
